# Apply updated simulation results (cases A-F run for the first time) to row 2
# of worksheet A-02.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -6239.788187305794
$ws.Range("C2").Value = 13134.93975530513
$ws.Range("D2").Value = -6239.788183591361
$ws.Range("E2").Value = -655.3633844079632
$ws.Range("F2").Value = 61.4919723795432
$ws.Range("G2").Value = 45.00550696016097
$ws.Range("H2").Value = 62.01355099356948
$ws.Range("I2").Value = 44.66130091520654
$ws.Range("J2").Value = 61.49197237924233
$ws.Range("K2").Value = 45.00550696011436
$ws.Range("L2").Value = 53.91398312472256
$ws.Range("O2").Value = 46.2358003495641
$ws.Range("P2").Value = 53.91398311941623
$ws.Range("R2").Value = 5.406849132284513
$ws.Range("S2").Value = -10.81369826143384
$ws.Range("T2").Value = 5.406849129149325
$ws.Range("X2").Value = -59.65924703725118
$ws.Range("Y2").Value = -77.89318276919613
$ws.Range("Z2").Value = -59.6592470372512
$ws.Range("AE2").Value = -5.406849132284513
$ws.Range("AF2").Value = 5.406849129149325
$ws.Range("AG2").Value = 5.406849132284513
$ws.Range("AH2").Value = -10.81369826143384
$ws.Range("AI2").Value = 5.406849129149325
$ws.Range("AJ2").Value = 5.406849132284513
$ws.Range("AK2").Value = -5.406849129149325
$ws.Range("AL2").Value = 9.116967865972464
$ws.Range("AM2").Value = -9.116967865972462
$ws.Range("AN2").Value = -59.65924703725118
$ws.Range("AO2").Value = -77.89318276919613
$ws.Range("AP2").Value = -59.6592470372512
$ws.Range("AQ2").Value = -9.116967865972464
$ws.Range("AR2").Value = 9.116967865972462
$ws.Range("AS2").Value = 61.4919723795432
$ws.Range("AT2").Value = 61.4919723795432
$ws.Range("AU2").Value = 62.01355099356937
$ws.Range("AV2").Value = 62.01355099356937
$ws.Range("AW2").Value = 62.01355099356948
$ws.Range("AX2").Value = 61.49197237924233
$ws.Range("AY2").Value = 61.49197237924233
$ws.Range("AZ2").Value = 45.00550696016097
$ws.Range("BA2").Value = 45.00550696016097
$ws.Range("BB2").Value = 44.66130091520654
$ws.Range("BC2").Value = 44.66130091532875
$ws.Range("BD2").Value = 44.6613009150841
$ws.Range("BE2").Value = 45.00550696011436
$ws.Range("BF2").Value = 45.00550696011442
$ws.Range("BG2").Value = 53.91398312472256
$ws.Range("BJ2").Value = 46.2358003495641
$ws.Range("BK2").Value = 53.91398311941623
